$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the H1 title "Play Fat Santa Slot for Free - Enjoy
#    Festive Gameplay".
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "*Meta description*") {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Fat Santa Slot for Free - Enjoy
#    Festive Gameplay" right before the final paragraph of the
#    document (the one that currently holds the "Create a feature
#    image..." image-prompt text).
#
#    We splice raw WordprocessingML in using InsertXML at a position
#    just before the last character of the preceding paragraph; this
#    produces a clean new <w:p> (leading empty run + single bold run)
#    without bleeding surrounding character formatting into it, and
#    without disturbing the paragraph we split from.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$precedingPara = $d.Paragraphs.Item($count - 1)
$splitPos = $precedingPara.Range.End - 1
$insertionRange = $d.Range($splitPos, $splitPos)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fat Santa Slot for Free - Enjoy Festive Gameplay</w:t></w:r></w:p>'
$insertionRange.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 3) Swap out the image-generation prompt text on the (now) final
#    paragraph for the meta-description copy, keeping its existing
#    italic run formatting untouched.
# ------------------------------------------------------------------
$oldText = "Create a feature image for Fat Santa that showcases the festive theme of the game. The image should be in a cartoon style and prominently feature a happy Maya warrior with glasses. The warrior should be positioned in the center of the image, surrounded by snow and Christmas decorations. Santa and his sleigh should be flying above the warrior, dropping cakes down onto the reels. The reels should also be visible in the image, displaying the various Christmas-themed symbols. Overall, the image should be fun, colorful, and help to convey the festive atmosphere of the game."
$newText = "Read our review of Fat Santa slots and play for free. Enjoy the festive theme, bonus features, and mobile compatibility of this cheerful game."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newText, 2)
